$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: change coefficients in the objective-value formula (z = 5x+4y instead of 11x+3y)
$ws.Range("L13").Formula = "=5*L11+4*L12"

# Row 15: remove the stray "B" label in F15 (point label text moves to E15 only)
$ws.Range("F15").ClearContents()

# Row 16: second constraint swapped in for the corner-point calc (c1 -> c2),
# with its coefficients/RHS updated to match (x+2y<=6 instead of 6x+4y<=24)
$ws.Range("E16").Value = "c2"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2
$ws.Range("I16").Value = 6

# Row 18: same objective-value formula change as row 13
$ws.Range("L18").Formula = "=5*L16+4*L17"

# Row 23: same objective-value formula change again
$ws.Range("L23").Formula = "=5*L21+4*L22"

# Update the active selection to match the saved view state
$ws.Range("N25").Select()

$wb.Application.Calculate()
